# Refatorando o codigo: corrige os dados de teste na planilha "Cadastro"
# e ajusta a aba ativa / seleções das planilhas.

$wb = $excel.ActiveWorkbook

# Planilha "Cadastro": corrige sobrenome e email usados nos testes
$wsCadastro = $wb.Worksheets.Item("Cadastro")
$wsCadastro.Range("A2").Value = "Jefrey"
$wsCadastro.Range("B2").Value = "Sales"
$wsCadastro.Range("C2").Value = "antsa@email.com"

# Planilha "Lupa": reseta a seleção para A2
$wsLupa = $wb.Worksheets.Item("Lupa")
$wsLupa.Activate()
$wsLupa.Range("A2").Select() | Out-Null

# Ativa a planilha "Cadastro" (torna-se a aba selecionada) e posiciona a seleção em A2
$wsCadastro.Activate()
$wsCadastro.Range("A2").Select() | Out-Null
